$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$rows = @(1819)
foreach ($r in $rows) {
    $cell = $ws.Cells.Item($r, 3)
    $old = $cell.Value
    $cell.Value = $old / 3
}
